# Applies the "Short description of the project and user stories" edit:
#  1. Adds a first-line indent (720 twips = 36pt) to the "Безусловно, можно
#     стараться..." paragraph in the combat section.
#  2. Rewrites the intro paragraph of the "ПОЛЬЗОВАТЕЛЬСКИЕ ИСТОРИИ" section
#     with a new short description (and gives it the same 720-twip first
#     line indent), and rewrites four of the five user-story paragraphs
#     that follow with new text.

$d = $word.ActiveDocument

# --- 1. Indent the "Безусловно, ..." paragraph -----------------------------
$rng = $d.Content
$rng.Find.Execute("Безусловно, можно ст", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $rng.Paragraphs.First
$para.Format.FirstLineIndent = 36

# --- 2. Rewrite the project-description paragraph + give it the indent -----
$old1 = "Вы капитан космического корабля, используя ваш ум, находите оптимальное вариант развития событий при решении внутриигровых ситуаций, диалогов и выходите победителем из схваток с вражеским кораблём. Выживите в суровых условиях неизведанного космоса и найдите дорогу домой!"
$new1 = "Вы капитан космического корабля, используя ваш интелект, находите оптимальный вариант решения внутриигровых ситуаций, диалогов и тем самым выходите победителем из схваток с вражеским кораблём. Выживите в суровых условиях неизведанного космоса и найдите дорогу домой!"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

$rng2 = $d.Content
$rng2.Find.Execute("используя ваш интелект", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para2 = $rng2.Paragraphs.First
$para2.Format.FirstLineIndent = 36

# --- 3. Rewrite the four user-story paragraphs ------------------------------
$old2 = "Как капитан корабля, я хочу прокладывать маршрут нажимая на экран и тем самым достигать необходимых мест на карте."
$new2 = "Как капитан корабля, я хочу прокладывать маршрут касанием точки экране и тем самым достигать необходимых мест на карте."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

$old3 = "Как капитан корабля, я хочу улучшать свой экипаж и корабль, дабы упростить процесс игры."
$new3 = "Как капитан корабля, я хочу иметь простую систему улучшения своего экипажа и корабля, дабы упростить процесс игры в дальнейшем."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

$old4 = "Как капитан корабля, я хочу получать вознагрождения за выполнение игровых заданий."
$new4 = "Как капитан корабля, я хочу получать вознагрождения соответствующие уровню сложности игрового задания."
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# This paragraph is immediately followed by an empty run with identical
# run formatting (just <w:rtl val="0"/>), and it is the very end of the
# story. A plain Find/Replace there causes the engine to merge the edited
# run with that trailing empty run (since their formatting becomes
# identical), silently dropping it. Toggle Bold on before the replace (so
# the edited run's formatting differs from its neighbour and the merge is
# skipped) and back off afterwards (which cleanly drops the now-redundant
# <w:b/> element) to keep the empty run intact.
$old5 = "Как капитан корабля, я хочу иметь возможно альтернативным способом получать награду."
$new5 = "Как капитан корабля, я хочу иметь возможность, кроме выполнения квестов, альтернативным способом получать игровую валюту и опыт."
$rng5 = $d.Content
$rng5.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng5.Bold = 1
$rng5.Text = $new5
$rng5b = $d.Content
$rng5b.Find.Execute($new5, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng5b.Bold = 0

Write-Output "done"
